$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45901, 5594.48067261627, 4864.4470356076, 5472, 6203.576611, 0.0642905829720348),
    @(45902, 5601.17759246688, 5063.30188236521, 2952, 6210.033796, 113.339920245764),
    @(45903, 5600.41202886578, 5071.36441285989, 2952, 6209.210973, 113.673473208088),
    @(45904, 5628.9357609315, 5061.51121424303, 2952, 6257.484242, 114.08582063798),
    @(45905, 5628.93576093312, 4328.63006870248, 2952, 6257.484254, 83.5491067403899),
    @(45906, 1142.23741561831, 2163.38890053204, 2952, 2113.11063, 7.5942547880721),
    @(45907, 970.493609902267, 2101.58201365382, 2952, 1974.370109, 6.39410469798133),
    @(45908, 6110.60514061168, 4990.37301472486, 2952, 6751.167337, 111.622300463049),
    @(45909, 6110.60514061168, 5180.70127280373, 2952, 6751.167337, 119.552644549669),
    @(45910, 6110.60514061168, 5193.72877250248, 2952, 6751.167337, 120.095457037117),
    @(45911, 6110.60514061168, 5288.62513526743, 2952, 6751.167337, 124.049472152323),
    @(45912, 6110.60514061168, 4654.92674024802, 2952, 6751.167337, 97.6453723598477),
    @(45913, 1154.34388801228, 2273.18936525638, 2952, 2198.215749, 15.2108844268374),
    @(45914, 1022.76404167669, 2184.78236862696, 2952, 2053.079618, 10.9624143729277)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
